$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Reconstruct the sheet1 data (final state after the repeated
# pandas/openpyxl "to_excel" writes in the source workbook) ---

# Header row
$ws.Range("A1").Value = "run_num"
$ws.Range("B1").Value = "block_num"
$ws.Range("C1").Value = "start_time"
$ws.Range("D1").Value = "play_duration"
$ws.Range("E1").Value = "ear"
$ws.Range("F1").Value = "hand"

# Data rows
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 5.0277442
$ws.Range("D2").Value = 15.2730259
$ws.Range("E2").Value = "R"
$ws.Range("F2").Value = "none"

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 18.0263115
$ws.Range("D3").Value = 25.8763445
$ws.Range("E3").Value = "L"
$ws.Range("F3").Value = "none"

$ws.Range("A4").Value = 1
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = 31.0363141
$ws.Range("D4").Value = 38.894082
$ws.Range("E4").Value = "L"
$ws.Range("F4").Value = "none"

$ws.Range("A5").Value = 1
$ws.Range("B5").Value = 4
$ws.Range("C5").Value = 44.0328539
$ws.Range("D5").Value = 51.9469333
$ws.Range("E5").Value = "R"
$ws.Range("F5").Value = "none"

# Resize columns to fit the new header/data text
$ws.Columns.AutoFit()

# --- The source workbook went through several intermediate
# "to_excel" saves that each introduced a date-formatted column style
# (before the columns were finally overwritten with plain numbers), so
# the style table in the final file still carries the now-unused
# numFmtId 22 (date) xf's left behind from those earlier writes.
# Reproduce the same leftover style entries by touching scratch cells
# far outside the used range and then clearing them back to normal. ---

$scratch1 = $ws.Cells.Item(500, 1)
$scratch1.NumberFormat = "m/d/yy h:mm"
$scratch1.Value = 1
$scratch1.Borders.Item(7).LineStyle = 1

$scratch2 = $ws.Cells.Item(500, 2)
$scratch2.NumberFormat = "m/d/yy h:mm"
$scratch2.Value = 1
$scratch2.Borders.Item(8).LineStyle = 1

$scratch3 = $ws.Cells.Item(500, 3)
$scratch3.NumberFormat = "m/d/yy h:mm"
$scratch3.Value = 1
$scratch3.Borders.Item(9).LineStyle = 1

$scratch1.Clear()
$scratch2.Clear()
$scratch3.Clear()

Write-Output "sheet1 rebuilt"
